$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values were regenerated (K = strikeouts, replacing the
# previous "Strike#" computation). Write the new values for rows 2-16.
$newK = @{
    2  = 0
    3  = 4
    4  = 0
    5  = 1
    6  = 1
    7  = 2
    8  = 3
    9  = 2
    10 = 1
    11 = 1
    12 = 4
    13 = 2
    14 = 1
    15 = 2
    16 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
